$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New trade row (row 15) appended below the existing data.
$ws.Range("A15").Value = "V"

# Force the date-looking text to stay as literal text (matching the
# existing rows, which store dates as plain strings, not date serials).
$ws.Range("B15").NumberFormat = "@"
$ws.Range("B15").Value = "08/19/20"
$ws.Range("B15").Style = "Normal"

$ws.Range("C15").Value = "BUY"
$ws.Range("D15").Value = 12
$ws.Range("E15").Value = 32
$ws.Range("F15").Value = -384
$ws.Range("G15").Value = 32
$ws.Range("H15").Value = 384
$ws.Range("I15").Value = 12

# REALIZED_PROFIT is blank for this row, like most other rows in the
# sheet. A bare quote forces Excel to materialize the (empty) cell
# instead of silently skipping it, and resetting the style keeps it
# free of the "quote prefix" formatting that would otherwise stick.
$ws.Range("J15").Value = "'"
$ws.Range("J15").Style = "Normal"
